$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.303.48"
$ws.Range("E2").Value = "  +1.24%  "

$ws.Range("D3").Value = "1.624.11"
$ws.Range("E3").Value = "  +1.49%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").Value = "1.848.74"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").Value = "1.627.86"
$ws.Range("E13").Value = "  +1.68%  "

$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("E15").Value = "  +0.98%  "

$ws.Range("D16").Value = "26.307.09"
$ws.Range("E16").Value = "  +1.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.74%  "

$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("E21").Value = "  +1.45%  "

$ws.Range("E22").Value = "  +1.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("E24").Value = "  +7.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("E28").Value = "  +0.86%  "

$ws.Range("E29").Value = "  +1.78%  "

$ws.Range("E30").Value = "  +10.73%  "

$ws.Range("E31").Value = "  +0.55%  "

$ws.Range("E32").Value = "  +2.70%  "

$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.69%  "

$ws.Range("D36").Value = "1.172.40"
$ws.Range("E36").Value = "  +4.37%  "

$ws.Range("E37").Value = "  +0.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.814"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("E41").Value = "  +1.45%  "

$ws.Range("E42").Value = "  +0.94%  "

$ws.Range("E43").Value = "  +3.60%  "

$ws.Range("D44").Value = "1.759.18"
$ws.Range("E44").Value = "  +1.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("E46").Value = "  +15.31%  "

$ws.Range("E47").Value = "  +1.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "

$ws.Range("E49").Value = "  +1.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("E51").Value = "  -0.28%  "
